$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

$ws.Cells.Item($row, 1).Value = "16ZHA7"
$ws.Cells.Item($row, 2).Value = "Cuchilla de limpieza"
$ws.Cells.Item($row, 3).Value = "Ricoh Aficio 220 270 1015 1018 1113 2015 2016 2018 2020 MP1500 MP1600 MP1900 MP2000 MP2001 MP2501 MP2554 MP3054 MP3554 MP4054"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E42-D42)*G42"
$ws.Cells.Item($row, 9).Formula = "=D42*F42"
$ws.Cells.Item($row, 10).Value = 0
